$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" date column (C) for rows 2-7 from 45184 to 45185
$ws.Range("C2:C7").Value = 45185
